$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: change A17's number format to match the other date rows (YYYY-MM-DD HH:MM:SS)
$ws.Range("A17").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 18: new data row
$ws.Range("A18").Value = 44526
$ws.Range("A18").NumberFormat = "YYYY-MM-DD"
$ws.Range("B18").Value = 65265.85
